$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 219, pushing the existing rows 219-230 down to 221-232.
$ws.Rows("219:220").Insert()

# Row 219: new weekly entry (Primera quality)
$ws.Cells.Item(219, 1).Value = 1
$ws.Cells.Item(219, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(219, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(219, 4).Value = 44578
$ws.Cells.Item(219, 5).Value = 15
$ws.Cells.Item(219, 6).Value = 100114014
$ws.Cells.Item(219, 7).Value = "Betarraga"
$ws.Cells.Item(219, 8).Value = "Sin especificar"
$ws.Cells.Item(219, 9).Value = "Primera"
$ws.Cells.Item(219, 10).Value = 1200
$ws.Cells.Item(219, 11).Value = 350
$ws.Cells.Item(219, 12).Value = 400
$ws.Cells.Item(219, 13).Value = 375
$ws.Cells.Item(219, 14).Value = "`$/paquete 4 unidades"
$ws.Cells.Item(219, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(219, 16).Value = 94
$ws.Cells.Item(219, 17).Value = 4
$ws.Cells.Item(219, 18).Value = "Hortaliza"

# Row 220: new weekly entry (Segunda quality)
$ws.Cells.Item(220, 1).Value = 1
$ws.Cells.Item(220, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(220, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(220, 4).Value = 44578
$ws.Cells.Item(220, 5).Value = 15
$ws.Cells.Item(220, 6).Value = 100114014
$ws.Cells.Item(220, 7).Value = "Betarraga"
$ws.Cells.Item(220, 8).Value = "Sin especificar"
$ws.Cells.Item(220, 9).Value = "Segunda"
$ws.Cells.Item(220, 10).Value = 1200
$ws.Cells.Item(220, 11).Value = 350
$ws.Cells.Item(220, 12).Value = 400
$ws.Cells.Item(220, 13).Value = 375
$ws.Cells.Item(220, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(220, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(220, 16).Value = 75
$ws.Cells.Item(220, 17).Value = 5
$ws.Cells.Item(220, 18).Value = "Hortaliza"
